$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
}

# Row 2
$ws.Range('D2').Value = '67.281.44'
$ws.Range('E2').Value = '  +1.05%  '

# Row 3
$ws.Range('D3').Value = '3.948.10'
$ws.Range('E3').Value = '  +4.18%  '

# Row 4
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
Set-TextValue 'D5' '471.67'
$ws.Range('E5').Value = '  +8.73%  '

# Row 6
Set-TextValue 'D6' '145.77'
$ws.Range('E6').Value = '  +4.24%  '

# Row 7
Set-TextValue 'D7' '0.623'
$ws.Range('E7').Value = '  -0.15%  '

# Row 8
$ws.Range('E8').Value = '  -0.14%  '

# Row 9
Set-TextValue 'D9' '0.732'
$ws.Range('E9').Value = '  -0.60%  '

# Row 10
$ws.Range('E10').Value = '  +7.36%  '

# Row 11
Set-TextValue 'D11' '0.0000339'
$ws.Range('E11').Value = '  +6.93%  '

# Row 12
Set-TextValue 'D12' '43.40'
$ws.Range('E12').Value = '  +1.83%  '

# Row 13
$ws.Range('D13').Value = '4.550.41'
$ws.Range('E13').Value = '  +3.03%  '

# Row 14
Set-TextValue 'D14' '10.34'
$ws.Range('E14').Value = '  -1.15%  '

# Row 15
Set-TextValue 'D15' '15.12'
$ws.Range('E15').Value = '  +0.29%  '

# Row 16
$ws.Range('D16').Value = '3.962.48'
$ws.Range('E16').Value = '  +4.71%  '

# Row 17
$ws.Range('E17').Value = '  -0.22%  '

# Row 18
Set-TextValue 'D18' '19.81'
$ws.Range('E18').Value = '  -0.75%  '

# Row 19
Set-TextValue 'D19' '1.16'
$ws.Range('E19').Value = '  +2.37%  '

# Row 20
$ws.Range('D20').Value = '67.483.38'
$ws.Range('E20').Value = '  +1.16%  '

# Row 21
Set-TextValue 'D21' '431.90'
$ws.Range('E21').Value = '  +5.06%  '

# Row 22
Set-TextValue 'D22' '3.38'
$ws.Range('E22').Value = '  +3.36%  '

# Row 23
Set-TextValue 'D23' '14.56'
$ws.Range('E23').Value = '  -0.80%  '

# Row 24
Set-TextValue 'D24' '87.21'
$ws.Range('E24').Value = '  +2.22%  '

# Row 25
Set-TextValue 'D25' '3.58'
$ws.Range('E25').Value = '  +7.52%  '

# Row 26
Set-TextValue 'D26' '38.64'
$ws.Range('E26').Value = '  +4.73%  '

# Row 27
Set-TextValue 'D27' '5.75'
$ws.Range('E27').Value = '  +2.43%  '

# Row 28
Set-TextValue 'D28' '10.21'
$ws.Range('E28').Value = '  +4.15%  '

# Row 29
Set-TextValue 'D29' '9.64'
$ws.Range('E29').Value = '  -0.66%  '

# Row 30
Set-TextValue 'D30' '729.93'
$ws.Range('E30').Value = '  +1.72%  '

# Row 31
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D31' '13.58'
$ws.Range('E31').Value = '  -2.09%  '

# Row 32
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D32' '0.132'
$ws.Range('E32').Value = '  -2.79%  '

# Row 33
$ws.Range('E33').Value = '  +1.27%  '

# Row 34
Set-TextValue 'D34' '42.83'
$ws.Range('E34').Value = '  +3.20%  '

# Row 35
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D35' '0.153'
$ws.Range('E35').Value = '  +2.16%  '

# Row 36
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D36' '57.89'
$ws.Range('E36').Value = '  +3.87%  '

# Row 37
$ws.Range('E37').Value = '  -0.02%  '

# Row 38
$ws.Range('D38').Value = '0.0₃0785'
$ws.Range('E38').Value = '  +12.96%  '

# Row 39
Set-TextValue 'D39' '5.39'
$ws.Range('E39').Value = '  -4.77%  '

# Row 40
$ws.Range('E40').Value = '  +0.73%  '

# Row 41
Set-TextValue 'D41' '3.05'
$ws.Range('E41').Value = '  +3.09%  '

# Row 42
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D42' '0.141'
$ws.Range('E42').Value = '  -0.85%  '

# Row 43
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D43' '2.58'
$ws.Range('E43').Value = '  -6.42%  '

# Row 44
$ws.Range('E44').Value = '  -0.20%  '

# Row 45
$ws.Range('E45').Value = '  +3.78%  '

# Row 46
$ws.Range('E46').Value = '  +4.34%  '

# Row 47
Set-TextValue 'D47' '2.19'
$ws.Range('E47').Value = '  +5.20%  '

# Row 48
Set-TextValue 'D48' '3.44'
$ws.Range('E48').Value = '  +3.17%  '

# Row 49
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D49' '3.18'
$ws.Range('E49').Value = '  -1.22%  '

# Row 50
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D50' '146.40'
$ws.Range('E50').Value = '  +2.98%  '

# Row 51
$ws.Range('E51').Value = '  +1.77%  '
